$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G (WORK_LOCATION) for the new ORGANISATION field.
# This shifts the existing G:P columns (WORK_LOCATION..RECRUITER) one to the
# right, to H:Q, and leaves G1/G2 blank for us to fill in.
$ws.Columns("G:G").Insert()

# Header row
$ws.Range("G1").Value = "ORGANISATION"

# Row 2 data updates
$ws.Range("C2").Value = "01_Reject-No Show"
$ws.Range("D2").Value = "Vinod Uttam Chavan"
$ws.Range("E2").Value = '["8208909261",""]'
$ws.Range("F2").Value = "vinodchavan3371@gmail.com"
$ws.Range("G2").Value = "Facile info serv Pvt"
$ws.Range("H2").Value = "Mumbai"
$ws.Range("I2").Value = 4.4
$ws.Range("J2").Value = 5.3
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = "2022-11-18 11:37:36"
$ws.Range("N2").Value = "2022-11-19 12:30:00"
